$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 9) with form-submission style data, matching
# the pattern of the previous entry row (row 8).
$ws.Range("A9").Value = "Nara"
$ws.Range("B9").Value = "simhan@sim.com"
$ws.Range("C9").Value = "nothing"
$ws.Range("D9").Value = "2025-10-02T08:47:45.794Z"
